$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# A1 ("IN") is unchanged text-wise; F1's label changes from "IsNumber" to "Comment"
# now that the per-row Y/IsNumber flags are going away (see below).
$ws.Range("F1").Value = "Comment"

# --- Extend the ROW()-1 helper formula down through row 64 ---
$ws.Range("A32:A64").Formula = "=ROW()-1"

# --- Drop the old per-row "Y" (IsNumber) flags in column F for rows 11-18, 21, 22, 24 ---
$ws.Range("F11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("F21").ClearContents()
$ws.Range("F22").ClearContents()
$ws.Range("F24").ClearContents()

# --- New CountCells puzzle rows 32-41 ---
$ws.Range("B32:E32").Value = 69
$ws.Range("F32").Value = "for CC"
$ws.Range("B33:E33").Value = 68
$ws.Range("F33").Value = "for CC"
$ws.Range("B34:E34").Value = 64
$ws.Range("F34").Value = "for CC"
$ws.Range("B35:E35").Value = 62
$ws.Range("F35").Value = "for CC"
$ws.Range("B36:E36").Value = 66
$ws.Range("F36").Value = "for CC"
$ws.Range("B37:E37").Value = 60
$ws.Range("F37").Value = "for CC"
$ws.Range("B38:E38").Value = 55
$ws.Range("F38").Value = "for CC"
$ws.Range("B39:E39").Value = 53
$ws.Range("F39").Value = "for CC"
$ws.Range("B40:E40").Value = 46
$ws.Range("F40").Value = "for CC"
$ws.Range("B41:E41").Value = 54
$ws.Range("F41").Value = "for CC"

# --- New animal-name puzzle rows 42-64 (EN in B & D, ES in C & E) ---
$ws.Range("B42").Value = "Buffalo"
$ws.Range("C42").Value = "Búfalo"
$ws.Range("D42").Value = "Buffalo"
$ws.Range("E42").Value = "Búfalo"
$ws.Range("B43").Value = "Camel"
$ws.Range("C43").Value = "Camello"
$ws.Range("D43").Value = "Camel"
$ws.Range("E43").Value = "Camello"
$ws.Range("B44").Value = "Deer"
$ws.Range("C44").Value = "Venado"
$ws.Range("D44").Value = "Deer"
$ws.Range("E44").Value = "Venado"
$ws.Range("B45").Value = "Elephant"
$ws.Range("C45").Value = "Elefante"
$ws.Range("D45").Value = "Elephant"
$ws.Range("E45").Value = "Elefante"
$ws.Range("B46").Value = "Giraffe"
$ws.Range("C46").Value = "Jirafa"
$ws.Range("D46").Value = "Giraffe"
$ws.Range("E46").Value = "Jirafa"
$ws.Range("B47").Value = "Gorilla"
$ws.Range("C47").Value = "Gorila"
$ws.Range("D47").Value = "Gorilla"
$ws.Range("E47").Value = "Gorila"
$ws.Range("B48").Value = "Horse"
$ws.Range("C48").Value = "Caballo"
$ws.Range("D48").Value = "Horse"
$ws.Range("E48").Value = "Caballo"
$ws.Range("B49").Value = "Kangaroo"
$ws.Range("C49").Value = "Canguro"
$ws.Range("D49").Value = "Kangaroo"
$ws.Range("E49").Value = "Canguro"
$ws.Range("B50").Value = "Leopard"
$ws.Range("C50").Value = "Leopardo"
$ws.Range("D50").Value = "Leopard"
$ws.Range("E50").Value = "Leopardo"
$ws.Range("B51").Value = "Rabbit"
$ws.Range("C51").Value = "Conejo"
$ws.Range("D51").Value = "Rabbit"
$ws.Range("E51").Value = "Conejo"
$ws.Range("B52").Value = "Squirrel"
$ws.Range("C52").Value = "Ardilla"
$ws.Range("D52").Value = "Squirrel"
$ws.Range("E52").Value = "Ardilla"
$ws.Range("B53").Value = "Whale"
$ws.Range("C53").Value = "Bellena"
$ws.Range("D53").Value = "Whale"
$ws.Range("E53").Value = "Bellena"
$ws.Range("B54").Value = "Dolphin"
$ws.Range("C54").Value = "Delfín"
$ws.Range("D54").Value = "Dolphin"
$ws.Range("E54").Value = "Delfín"
$ws.Range("B55").Value = "Crow"
$ws.Range("C55").Value = "Cuervo"
$ws.Range("D55").Value = "Crow"
$ws.Range("E55").Value = "Cuervo"
$ws.Range("B56").Value = "Dove"
$ws.Range("C56").Value = "Paloma"
$ws.Range("D56").Value = "Dove"
$ws.Range("E56").Value = "Paloma"
$ws.Range("B57").Value = "Eagle"
$ws.Range("C57").Value = "Águila"
$ws.Range("D57").Value = "Eagle"
$ws.Range("E57").Value = "Águila"
$ws.Range("B58").Value = "Flamingo"
$ws.Range("C58").Value = "Flamenco"
$ws.Range("D58").Value = "Flamingo"
$ws.Range("E58").Value = "Flamenco"
$ws.Range("B59").Value = "Hummingbird"
$ws.Range("C59").Value = "Colibrí"
$ws.Range("D59").Value = "Hummingbird"
$ws.Range("E59").Value = "Colibrí"
$ws.Range("B60").Value = "Ostrich"
$ws.Range("C60").Value = "Avestruz"
$ws.Range("D60").Value = "Ostrich"
$ws.Range("E60").Value = "Avestruz"
$ws.Range("B61").Value = "Parrot"
$ws.Range("C61").Value = "Papagayo"
$ws.Range("D61").Value = "Parrot"
$ws.Range("E61").Value = "Papagayo"
$ws.Range("B62").Value = "Pelican"
$ws.Range("C62").Value = "Pelicano"
$ws.Range("D62").Value = "Pelican"
$ws.Range("E62").Value = "Pelicano"
$ws.Range("B63").Value = "Pigeon"
$ws.Range("C63").Value = "Paloma"
$ws.Range("D63").Value = "Pigeon"
$ws.Range("E63").Value = "Paloma"
$ws.Range("B64").Value = "Quetzal"
$ws.Range("C64").Value = "Quetzal"
$ws.Range("D64").Value = "Quetzal"
$ws.Range("E64").Value = "Quetzal"

# --- Update the saved selection / scroll position to match the edited view ---
$ws.Range("A36:A64").Select()

Write-Host "edit applied"
